$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("Data")
$wsCode = $wb.Worksheets.Item("Codebook")

# --- Codebook sheet: fix typo for the "Age" definition ---
$wsCode.Range("B6").Value = "age in number"

# --- Data sheet: fill in the two new variables (Smoking answers, Age numbers) ---
$smoking = @{
    2  = "Y"
    3  = "N"
    4  = "Y"
    5  = "Y"
    6  = "Y"
    7  = "NA"
    8  = "N"
    9  = "N"
    10 = "N"
    11 = "Y"
    12 = "N"
    13 = "NA"
    14 = "N"
    15 = "Y"
}

$age = @{
    2  = 76
    3  = 11
    4  = 44
    5  = 48
    6  = 79
    7  = 45
    8  = 30
    9  = 38
    10 = 27
    11 = 19
    12 = 51
    13 = 20
    14 = 35
    15 = 14
}

foreach ($row in 2..15) {
    $wsData.Range("D$row").Value = $smoking[$row]
    $wsData.Range("E$row").Value = $age[$row]
}

# --- Codebook sheet: broaden allowed values for Smoking now that "NA" shows up ---
$wsCode.Range("C5").Value = "Y/N/NA"

# --- View state: Data tab now active, Codebook tab no longer the active one ---
[void]$wsCode.Range("C5").Select()
[void]$wsData.Select()
[void]$wsData.Range("I18").Select()
